# Swap the occurrence data between row 2 and row 3.
# (Id, Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn, Auktor, Ost, Nord)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ between row 2 and row 3 and must be swapped: A, B, E, F, G, H, Q, R
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $tmp = $cell2.Value2
    $cell2.Value = $cell3.Value2
    $cell3.Value = $tmp
}
